$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new review row (row 11) right below the current last row (row 10)
$ws.Range("A11").Value = "com.singleton.strechy"
$ws.Range("B11").Value = "stretchy"
$ws.Range("C11").Value = "sixsevensix67676@gmail.com"
$ws.Range("D11").Value = "stevewonder3001@gmail.com"
$ws.Range("E11").Value = "27/5/2019 15:59"
$ws.Range("F11").Value = "This store is full with a lot of good games, but this cars game is the best. Guaranteed!!"

# Add the mailto hyperlinks for the two email cells in the new row
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:sixsevensix67676@gmail.com", "", "", "sixsevensix67676@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com") | Out-Null

# Adding hyperlinks auto-applies Excel's built-in "Hyperlink" cell style; restore
# the new row's look to match the rest of the table (same formatting as row 10)
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)

# Clean up the now-unused "Hyperlink" named style that Excel auto-created
try {
    $wb.Styles.Item("Hyperlink").Delete() | Out-Null
} catch {
}

# Keep the active selection in sync with the new last cell
$ws.Range("F11").Select() | Out-Null
